$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.046891843372023
$ws.Range("D2").Value = 1.051432107831127
$ws.Range("E2").Value = 1.044365434243451
$ws.Range("F2").Value = 1.06160513625515
$ws.Range("I2").Value = 1.042651040069261
$ws.Range("J2").Value = 1.051943820030507
$ws.Range("K2").Value = 1.054183626539406
$ws.Range("L2").Value = 1.047136688755698
$ws.Range("M2").Value = 1.064328743368869
$ws.Range("N2").Value = 1.053437700914593
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.048414639353447
$ws.Range("D3").Value = 1.052589044441482
$ws.Range("E3").Value = 1.04568232934531
$ws.Range("F3").Value = 1.062896389572507
$ws.Range("I3").Value = 1.043056240418593
$ws.Range("J3").Value = 1.053112047231978
$ws.Range("K3").Value = 1.055152322889519
$ws.Range("L3").Value = 1.048263487955579
$ws.Range("M3").Value = 1.065433450052586
$ws.Range("N3").Value = 1.054607587132687
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.049398523904167
$ws.Range("D4").Value = 1.053336191011605
$ws.Range("E4").Value = 1.046533292221892
$ws.Range("F4").Value = 1.063730733772709
$ws.Range("I4").Value = 1.043316155072624
$ws.Range("J4").Value = 1.053866121976085
$ws.Range("K4").Value = 1.055777090590245
$ws.Range("L4").Value = 1.048990907936311
$ws.Range("M4").Value = 1.066146535464569
$ws.Range("N4").Value = 1.055362732749425
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.049811805674939
$ws.Range("D5").Value = 1.053649944607727
$ws.Range("E5").Value = 1.046890766041617
$ws.Range("F5").Value = 1.064081214002134
$ws.Range("I5").Value = 1.043424880257166
$ws.Range("J5").Value = 1.054182698612457
$ws.Range("K5").Value = 1.056039258198641
$ws.Range("L5").Value = 1.049296315009572
$ws.Range("M5").Value = 1.066445906270923
$ws.Range("N5").Value = 1.055679758960933
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.049881177595821
$ws.Range("D6").Value = 1.053702605018349
$ws.Range("E6").Value = 1.046950771748266
$ws.Range("F6").Value = 1.06414004501755
$ws.Range("I6").Value = 1.043443103919573
$ws.Range("J6").Value = 1.054235827796045
$ws.Range("K6").Value = 1.056083248999237
$ws.Range("L6").Value = 1.049347570883427
$ws.Range("M6").Value = 1.066496148023668
$ws.Range("N6").Value = 1.055732963594055
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.049404047534764
$ws.Range("D7").Value = 1.05334038475936
$ws.Range("E7").Value = 1.046538069862076
$ws.Range("F7").Value = 1.063735417992968
$ws.Range("I7").Value = 1.043317609994533
$ws.Range("J7").Value = 1.053870353796338
$ws.Range("K7").Value = 1.055780595585775
$ws.Range("L7").Value = 1.048994990366798
$ws.Range("M7").Value = 1.066150537278291
$ws.Range("N7").Value = 1.055366970579348
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047406787673576
$ws.Range("D8").Value = 1.051823406998677
$ws.Range("E8").Value = 1.044810728061773
$ws.Range("F8").Value = 1.06204176942545
$ws.Range("I8").Value = 1.042788452234851
$ws.Range("J8").Value = 1.052339013590752
$ws.Range("K8").Value = 1.054511427874758
$ws.Range("L8").Value = 1.047517849907632
$ws.Range("M8").Value = 1.064702446001549
$ws.Range("N8").Value = 1.053833455695023
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043875783106667
$ws.Range("D9").Value = 1.049138832726835
$ws.Range("E9").Value = 1.041757812077018
$ws.Range("F9").Value = 1.059048051633668
$ws.Range("I9").Value = 1.041838469783697
$ws.Range("J9").Value = 1.04962620178165
$ws.Range("K9").Value = 1.052259137333972
$ws.Range("L9").Value = 1.044901720827838
$ws.Range("M9").Value = 1.062137232399849
$ws.Range("N9").Value = 1.051116791381988
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.041513507968931
$ws.Range("D10").Value = 1.04734110220531
$ws.Range("E10").Value = 1.039716033609585
$ws.Range("F10").Value = 1.057045664773425
$ws.Range("I10").Value = 1.041193215497435
$ws.Range("J10").Value = 1.047807625382248
$ws.Range("K10").Value = 1.050746663781903
$ws.Range("L10").Value = 1.043148407562621
$ws.Range("M10").Value = 1.060417717262961
$ws.Range("N10").Value = 1.049295632395505
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.040488545539752
$ws.Range("D11").Value = 1.046560695358697
$ws.Range("E11").Value = 1.038830296018702
$ws.Range("F11").Value = 1.056176976006995
$ws.Range("I11").Value = 1.040910951204822
$ws.Range("J11").Value = 1.047017701395699
$ws.Range("K11").Value = 1.050089085163962
$ws.Range("L11").Value = 1.042386938810727
$ws.Range("M11").Value = 1.059670857214731
$ws.Range("N11").Value = 1.04850458662629
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.040107506618314
$ws.Range("D12").Value = 1.046270514361233
$ws.Range("E12").Value = 1.038501040911495
$ws.Range("F12").Value = 1.055854053177732
$ws.Range("I12").Value = 1.040805672438857
$ws.Range("J12").Value = 1.046723910613474
$ws.Range("K12").Value = 1.04984442475844
$ws.Range("L12").Value = 1.042103747341903
$ws.Range("M12").Value = 1.059393088841188
$ws.Range("N12").Value = 1.048210378627456
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.040189255495487
$ws.Range("D13").Value = 1.046332772985565
$ws.Range("E13").Value = 1.038571678802892
$ws.Range("F13").Value = 1.05592333278012
$ws.Range("I13").Value = 1.040828274742053
$ws.Range("J13").Value = 1.046786946944579
$ws.Range("K13").Value = 1.049896923728951
$ws.Range("L13").Value = 1.042164508731452
$ws.Range("M13").Value = 1.059452687093292
$ws.Range("N13").Value = 1.048273504477382
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.040457055372342
$ws.Range("D14").Value = 1.046536715126069
$ws.Range("E14").Value = 1.038803084886449
$ws.Range("F14").Value = 1.056150288306709
$ws.Range("I14").Value = 1.040902257683535
$ws.Range("J14").Value = 1.046993424306506
$ws.Range("K14").Value = 1.050068869782402
$ws.Range("L14").Value = 1.042363537254727
$ws.Range("M14").Value = 1.059647903995302
$ws.Range("N14").Value = 1.048480275060845
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04062201273021
$ws.Range("D15").Value = 1.04666233025872
$ws.Range("E15").Value = 1.038945628105598
$ws.Range("F15").Value = 1.056290089359464
$ws.Range("I15").Value = 1.040947783522324
$ws.Range("J15").Value = 1.047120591545899
$ws.Range("K15").Value = 1.050174757395258
$ws.Range("L15").Value = 1.042486118974312
$ws.Range("M15").Value = 1.059768136842312
$ws.Range("N15").Value = 1.048607622892305
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.041581486682528
$ws.Range("D16").Value = 1.047392853057272
$ws.Range("E16").Value = 1.039774782052985
$ws.Range("F16").Value = 1.05710328164112
$ws.Range("I16").Value = 1.041211887848343
$ws.Range("J16").Value = 1.047859997433672
$ws.Range("K16").Value = 1.050790248424642
$ws.Range("L16").Value = 1.043198895289915
$ws.Range("M16").Value = 1.060467234967928
$ws.Range("N16").Value = 1.04934807882125
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042182775337301
$ws.Range("D17").Value = 1.047850556993545
$ws.Range("E17").Value = 1.040294445850455
$ws.Range("F17").Value = 1.057612931520904
$ws.Range("I17").Value = 1.041376784511287
$ws.Range("J17").Value = 1.048323141430936
$ws.Range("K17").Value = 1.0511756115209
$ws.Range("L17").Value = 1.043645387941758
$ws.Range("M17").Value = 1.060905141152969
$ws.Range("N17").Value = 1.049811880536113
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042533296627153
$ws.Range("D18").Value = 1.048117337549991
$ws.Range("E18").Value = 1.040597399782303
$ws.Range("F18").Value = 1.057910043678906
$ws.Range("I18").Value = 1.041472689712921
$ws.Range("J18").Value = 1.048593047799275
$ws.Range("K18").Value = 1.051400130296418
$ws.Range("L18").Value = 1.04390560073675
$ws.Range("M18").Value = 1.061160343041236
$ws.Range("N18").Value = 1.050082170202454
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042652781655352
$ws.Range("D19").Value = 1.048208270710544
$ws.Range("E19").Value = 1.040700672803077
$ws.Range("F19").Value = 1.058011324641948
$ws.Range("I19").Value = 1.041505344139166
$ws.Range("J19").Value = 1.048685038790039
$ws.Range("K19").Value = 1.051476641975368
$ws.Range("L19").Value = 1.043994289660612
$ws.Range("M19").Value = 1.061247322919107
$ws.Range("N19").Value = 1.050174291830976
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042118283532523
$ws.Range("D20").Value = 1.047801469431117
$ws.Range("E20").Value = 1.040238707152225
$ws.Range("F20").Value = 1.057558267300921
$ws.Range("I20").Value = 1.041359121252323
$ws.Range("J20").Value = 1.048273475099207
$ws.Range("K20").Value = 1.051134292331282
$ws.Range("L20").Value = 1.043597506187009
$ws.Range("M20").Value = 1.060858180905226
$ws.Range("N20").Value = 1.049762143672496
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.040378203955252
$ws.Range("D21").Value = 1.046476667665801
$ws.Range("E21").Value = 1.038734948597391
$ws.Range("F21").Value = 1.056083462600866
$ws.Range("I21").Value = 1.040880483533537
$ws.Range("J21").Value = 1.04693263235646
$ws.Range("K21").Value = 1.050018247220348
$ws.Range("L21").Value = 1.042304937977349
$ws.Range("M21").Value = 1.059590427228905
$ws.Range("N21").Value = 1.048419396779257
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.039282277736448
$ws.Range("D22").Value = 1.04564195527727
$ws.Range("E22").Value = 1.037788010077768
$ws.Range("F22").Value = 1.055154726511083
$ws.Range("I22").Value = 1.04057703699608
$ws.Range("J22").Value = 1.046087401723353
$ws.Range("K22").Value = 1.049314191060404
$ws.Range("L22").Value = 1.041490231886776
$ws.Range("M22").Value = 1.058791304689992
$ws.Range("N22").Value = 1.047572965821702
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.039863429317034
$ws.Range("D23").Value = 1.046084620669507
$ws.Range("E23").Value = 1.038290141476147
$ws.Range("F23").Value = 1.055647208384163
$ws.Range("I23").Value = 1.040738138423937
$ws.Range("J23").Value = 1.046535684369782
$ws.Range("K23").Value = 1.049687649612008
$ws.Range("L23").Value = 1.041922316435048
$ws.Range("M23").Value = 1.059215129608161
$ws.Range("N23").Value = 1.048021885080908
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.042147425221936
$ws.Range("D24").Value = 1.04782365057624
$ws.Range("E24").Value = 1.04026389355661
$ws.Range("F24").Value = 1.057582968197432
$ws.Range("I24").Value = 1.0413671033724
$ws.Range("J24").Value = 1.048295917911318
$ws.Range("K24").Value = 1.051152963488832
$ws.Range("L24").Value = 1.043619142567273
$ws.Range("M24").Value = 1.06087940090411
$ws.Range("N24").Value = 1.049784618355975
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.044790054452863
$ws.Range("D25").Value = 1.049834249553238
$ws.Range("E25").Value = 1.04254818261813
$ws.Range("F25").Value = 1.05982313417959
$ws.Range("I25").Value = 1.04208615584652
$ws.Range("J25").Value = 1.050329271093771
$ws.Range("K25").Value = 1.052843316668414
$ws.Range("L25").Value = 1.045579654215102
$ws.Range("M25").Value = 1.062802031343123
$ws.Range("N25").Value = 1.051820859133176
